$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.857 (0.855 ± 0.002)"
$ws.Range("C2").Value = "00:00:58 (00:01:25 ± 00:00:21)"
$ws.Range("D2").Value = "00:00:09 (00:00:09 ± 00:00:00)"
$ws.Range("B3").Value = "0.981 (0.967 ± 0.009)"
$ws.Range("C3").Value = "00:00:07 (00:00:08 ± 00:00:01)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B4").Value = "0.969 (0.945 ± 0.015)"
$ws.Range("C4").Value = "00:00:24 (00:00:37 ± 00:00:09)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B5").Value = "0.989 (0.972 ± 0.007)"
$ws.Range("C5").Value = "00:05:13 (00:05:22 ± 00:00:07)"
$ws.Range("D5").Value = "00:00:00 (00:00:02 ± 00:00:01)"
$ws.Range("B6").Value = "0.915 (0.881 ± 0.023)"
$ws.Range("C6").Value = "00:04:58 (00:05:02 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:00)"
$ws.Range("B7").Value = "0.954 (0.947 ± 0.009)"
$ws.Range("C7").Value = "00:05:00 (00:05:01 ± 00:00:00)"
$ws.Range("D7").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B9").Value = "0.985 (0.400 ± 0.329)"
$ws.Range("C9").Value = "00:04:59 (00:05:00 ± 00:00:00)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B11").Value = "0.985 (0.968 ± 0.009)"
$ws.Range("C11").Value = "00:05:01 (00:05:03 ± 00:00:01)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B13").Value = "0.985 (0.966 ± 0.008)"
$ws.Range("C13").Value = "00:00:57 (00:01:08 ± 00:00:07)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B14").Value = "0.962 (0.944 ± 0.012)"
$ws.Range("C14").Value = "00:00:26 (00:00:28 ± 00:00:01)"
$ws.Range("D14").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B15").Value = "0.969 (0.956 ± 0.010)"
$ws.Range("C15").Value = "00:04:33 (00:04:51 ± 00:00:08)"
$ws.Range("D15").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B16").Value = "0.965 (0.947 ± 0.011)"
$ws.Range("C16").Value = "00:00:11 (00:00:12 ± 00:00:00)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 ± 00:00:00)"
